# This script updates the weekly price records (rows 2-23) of the
# "Albahaca" subset sheet. The rows were re-ordered (a different weekly
# grouping), so for each destination row we write the new Fecha (D),
# Volumen (J), Precio minimo (K), Precio maximo (L), Precio promedio
# ponderado (M) and Precio $/Kg (P) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, PrecioKg
$data = @(
    @(2,  44260, 250, 900,  1000, 950,  950),
    @(3,  44250, 250, 1000, 1200, 1100, 1100),
    @(4,  44362, 250, 2800, 3000, 2900, 2900),
    @(5,  44432, 300, 2300, 2500, 2400, 2400),
    @(6,  44349, 250, 2800, 3000, 2900, 2900),
    @(7,  44498, 270, 2000, 2300, 2150, 2150),
    @(8,  44249, 200, 900,  1000, 950,  950),
    @(9,  44365, 250, 2400, 2500, 2450, 2450),
    @(10, 44224, 200, 750,  800,  775,  775),
    @(11, 44435, 300, 2300, 2500, 2400, 2400),
    @(12, 44313, 250, 900,  1000, 950,  950),
    @(13, 44274, 250, 1000, 1200, 1100, 1100),
    @(14, 44417, 250, 4000, 4500, 4250, 4250),
    @(15, 44326, 200, 2700, 2800, 2750, 2750),
    @(16, 44474, 250, 2000, 2500, 2250, 2250),
    @(17, 44376, 270, 2400, 2500, 2437, 2437),
    @(18, 44280, 250, 1400, 1500, 1450, 1450),
    @(19, 44442, 240, 2300, 2500, 2400, 2400),
    @(20, 44292, 250, 1800, 2000, 1900, 1900),
    @(21, 44494, 200, 2400, 2500, 2450, 2450),
    @(22, 44330, 250, 2800, 3000, 2900, 2900),
    @(23, 44302, 200, 900,  1000, 950,  950)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 4).Value  = $row[1]   # D - Fecha
    $ws.Cells.Item($r, 10).Value = $row[2]   # J - Volumen
    $ws.Cells.Item($r, 11).Value = $row[3]   # K - Precio minimo
    $ws.Cells.Item($r, 12).Value = $row[4]   # L - Precio maximo
    $ws.Cells.Item($r, 13).Value = $row[5]   # M - Precio promedio ponderado
    $ws.Cells.Item($r, 16).Value = $row[6]   # P - Precio $/Kg
}
